# Parameter_HeatingTechnology_Cost.xlsx edit:
#   - capex calculation revised: lifetime is replaced by payback time
#   - Add a new "payback_time" column (=30) to Table1 on sheet "Tabelle1"
#   - Remove the obsolete helper sheet "Sheet1"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Expand the table by one column (A1:N76 -> A1:O76)
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:O76"))

# Name the new column header; this also updates the table column's name
$ws.Range("O1").Value = "payback_time"

# Fill the new column with the payback time value (30) for every data row
for ($r = 2; $r -le 76; $r++) {
    $ws.Cells.Item($r, 15).Value = 30
}

# Remove the now unused "Sheet1" worksheet
$helper = $wb.Worksheets.Item("Sheet1")
$helper.Delete()
